# Fruta / hortaliza, semanal
# Insert a new weekly record at row 126 (pushing the existing rows 126-151
# down to 127-152, a new trailing row 152 appears carrying what used to be
# row 151's data) and populate the newly inserted row with this week's
# figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 126:151 down to 127:152, carrying along their values/styles.
$ws.Rows(126).Insert()

# Populate the newly-opened row 126 with the new weekly entry.
$ws.Cells.Item(126, 1).Value = 11
$ws.Cells.Item(126, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(126, 3).Value = "Bíobío"
$ws.Cells.Item(126, 4).Value = 44637
$ws.Cells.Item(126, 5).Value = 8
$ws.Cells.Item(126, 6).Value = 100112003
$ws.Cells.Item(126, 7).Value = "Ajo"
$ws.Cells.Item(126, 8).Value = "Chino"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 150
$ws.Cells.Item(126, 11).Value = 17000
$ws.Cells.Item(126, 12).Value = 18000
$ws.Cells.Item(126, 13).Value = 17467
$ws.Cells.Item(126, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(126, 15).Value = "China"
$ws.Cells.Item(126, 16).Value = 1747
$ws.Cells.Item(126, 17).Value = 10
$ws.Cells.Item(126, 18).Value = "Hortaliza"
